$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> (new D value or $null if unchanged, new E value or $null if unchanged)
$updates = @{
    2  = @("26.787.09", "  -1.71%  ")
    3  = @("1.869.29", "  -1.95%  ")
    4  = @("1.001", "  -0.02%  ")
    5  = @("300.90", "  -2.19%  ")
    6  = @("1.000", "  -0.04%  ")
    7  = @($null, "  +1.65%  ")
    8  = @("0.3739", "  -1.87%  ")
    9  = @("0.07183", "  -1.41%  ")
    10 = @("21.63", "  -0.19%  ")
    11 = @("0.8884", "  -1.75%  ")
    12 = @("0.08155", "  -0.37%  ")
    13 = @("1.881.55", "  +25.25%  ")
    14 = @("92.80", "  -3.74%  ")
    15 = @("5.303", "  -1.21%  ")
    16 = @("1.002", "  +0.09%  ")
    17 = @("14.82", "  +0.48%  ")
    18 = @("0.000008500", "  -1.90%  ")
    19 = @($null, "  -0.02%  ")
    20 = @("26.813.55", $null)
    21 = @("4.987", "  -2.54%  ")
    22 = @($null, "  -1.94%  ")
    23 = @("6.384", "  -1.86%  ")
    24 = @("2.315", "  -1.22%  ")
    25 = @("145.93", "  -2.91%  ")
    26 = @("1.735", "  -0.23%  ")
    27 = @("18.00", "  -1.44%  ")
    28 = @("113.82", "  -2.48%  ")
    29 = @("4.718", "  -2.65%  ")
    30 = @("4.623", "  -4.69%  ")
    31 = @("0.09152", "  -1.06%  ")
    32 = @("0.8038", "  -3.83%  ")
    33 = @("0.05027", "  -0.84%  ")
    34 = @("1.172", "  -4.76%  ")
    35 = @("2.933", "  -1.72%  ")
    36 = @("0.6114", "  +5.93%  ")
    37 = @("2.702", "  -1.24%  ")
    38 = @("3.195", "  -4.79%  ")
    39 = @($null, "  -2.55%  ")
    40 = @("1.066", "  -1.22%  ")
    41 = @("6.517", "  -0.76%  ")
    42 = @("0.5207", "  +5.89%  ")
    43 = @("8.771", "  -3.58%  ")
    44 = @("114.73", "  -2.18%  ")
    45 = @("0.1491", "  -2.05%  ")
    46 = @("1.000", "  -0.07%  ")
    47 = @("1.653", "  +0.48%  ")
    48 = @("9.950", "  -2.32%  ")
    49 = @("37.61", "  -3.32%  ")
    50 = @("0.06053", "  -0.06%  ")
    51 = @("62.19", "  -3.53%  ")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $dVal = $vals[0]
    $eVal = $vals[1]

    if ($null -ne $dVal) {
        $cell = $ws.Cells.Item($row, 4)
        # Preserve the cell's existing look (e.g. "1.000" / "0.000008500")
        # by temporarily forcing text entry, then restore the original
        # cell style so no stray formatting is introduced.
        $origStyle = $cell.Style
        $cell.NumberFormat = "@"
        $cell.Value = $dVal
        $cell.Style = $origStyle
    }
    if ($null -ne $eVal) {
        $cell = $ws.Cells.Item($row, 5)
        $cell.Value = $eVal
    }
}
